$wb = $excel.ActiveWorkbook

# --- Sheet "compare_models": update column I (MAPE / TT (Sec)?) values ---
$wsCompare = $wb.Worksheets.Item("compare_models")

$wsCompare.Range("I2").Value  = 0.09
$wsCompare.Range("I3").Value  = 0.054
$wsCompare.Range("I4").Value  = 0.088
$wsCompare.Range("I5").Value  = 0.06
$wsCompare.Range("I6").Value  = 0.032
$wsCompare.Range("I7").Value  = 1.08
$wsCompare.Range("I9").Value  = 0.024
$wsCompare.Range("I11").Value = 0.03
$wsCompare.Range("I12").Value = 0.544
$wsCompare.Range("I13").Value = 0.022
$wsCompare.Range("I15").Value = 0.02
$wsCompare.Range("I18").Value = 0.02
$wsCompare.Range("I19").Value = 0.022

# --- Sheet "pred_final": update row 2, columns C..H ---
$wsPredFinal = $wb.Worksheets.Item("pred_final")

$wsPredFinal.Range("C2").Value = 1.3982
$wsPredFinal.Range("D2").Value = 5.2116
$wsPredFinal.Range("E2").Value = 2.2829
$wsPredFinal.Range("F2").Value = 0.9921
$wsPredFinal.Range("G2").Value = 0.0441
$wsPredFinal.Range("H2").Value = 0.0238
